# Update the K column (column G) values on the active sheet to reflect
# the regenerated save_data (K instead of Strike#, recalculated std/mean,
# and rewritten s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K (column G) value
$kValues = @{
    2  = 1
    3  = 5
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 0
    9  = 5
    10 = 2
    11 = 6
    12 = 2
    13 = 5
    14 = 3
    15 = 3
    16 = 2
    17 = 5
    18 = 4
    19 = 1
    20 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
